$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currencies")

$data = @(
    @('AED', 'د.إ'),
    @('AFS', 'Af'),
    @('AKZ', 'Kz'),
    @('ARS', '$'),
    @('AUD', '$'),
    @('AZN', 'ман.'),
    @('BAM', 'KM'),
    @('BDT', '৳'),
    @('BGN', 'лв.'),
    @('BHD', '.د.ب'),
    @('BOB', 'Bs.'),
    @('BRL', 'R$'),
    @('BTC', 'Ƀ'),
    @('CAD', 'C$'),
    @('CHF', 'Fr.'),
    @('CLP', '$'),
    @('CNY', '¥'),
    @('COP', '$'),
    @('CRC', '₡'),
    @('CZK', 'Kč'),
    @('DKK', 'kr'),
    @('DOP', 'RD$'),
    @('EGP', 'ج.م'),
    @('EUR', 'â‚¬'),
    @('GBP', '£'),
    @('GTQ', 'Q'),
    @('HKD', '$'),
    @('HNL', 'L'),
    @('HRK', 'kn'),
    @('HUF', 'Ft'),
    @('IDR', 'Rp'),
    @('ILS', '₪'),
    @('INR', '₹'),
    @('JPY', '¥'),
    @('KES', 'KSh'),
    @('KRW', '₩'),
    @('KWD', 'د.ك'),
    @('KZT', 'T'),
    @('LBP', 'ل.ل.'),
    @('LKR', 'Rs'),
    @('MAD', 'درهم'),
    @('MXN', '$'),
    @('MYR', 'RM'),
    @('NGN', '₦'),
    @('NIO', 'C$'),
    @('NOK', 'kr'),
    @('NZD', '$'),
    @('OMR', 'ر.ع.'),
    @('PAB', '$'),
    @('PEN', 'S/.'),
    @('PHP', '₱'),
    @('PKR', 'Rs'),
    @('PLN', 'zł'),
    @('PYG', 'Gs.'),
    @('RON', 'lei'),
    @('RUB', '₽'),
    @('RWF', 'R₣'),
    @('SAR', 'ر.س'),
    @('SEK', 'kr'),
    @('SGD', '$'),
    @('SVC', '₡'),
    @('THB', '฿'),
    @('TND', 'د.ت'),
    @('TRY', '₺'),
    @('TWD', '$'),
    @('TZS', 'Tsh'),
    @('UAH', '₴'),
    @('USD', '$'),
    @('UYU', '$U'),
    @('VEF', 'Bs.'),
    @('VND', '₫'),
    @('XOF', 'CFA'),
    @('ZAR', 'R')
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Rows.Item($r).RowHeight = 25
}
